$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (data rows only) to text format so numeric-looking values (e.g. "6.37")
# stay as text, matching the workbook's existing inlineStr/text storage for that column.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '69.141.48'
$ws.Range("E2").Value = '  +0.59%  '
$ws.Range("D3").Value = '3.759.25'
$ws.Range("E3").Value = '  +0.15%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '602.75'
$ws.Range("E5").Value = '  +0.09%  '
$ws.Range("D6").Value = '166.56'
$ws.Range("E6").Value = '  -1.56%  '
$ws.Range("D7").Value = '3.757.92'
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E10").Value = '  +4.09%  '
$ws.Range("D11").Value = '6.37'
$ws.Range("E11").Value = '  +0.83%  '
$ws.Range("D12").Value = '0.460'
$ws.Range("E12").Value = '  -0.84%  '
$ws.Range("D13").Value = '37.72'
$ws.Range("E13").Value = '  -1.57%  '
$ws.Range("E14").Value = '  +0.40%  '
$ws.Range("D15").Value = '4.388.06'
$ws.Range("E15").Value = '  +0.15%  '
$ws.Range("D16").Value = '3.758.75'
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("D17").Value = '69.131.73'
$ws.Range("E17").Value = '  +0.56%  '
$ws.Range("D18").Value = '7.41'
$ws.Range("E18").Value = '  +1.57%  '
$ws.Range("D19").Value = '17.78'
$ws.Range("E19").Value = '  +3.55%  '
$ws.Range("E20").Value = '  -0.82%  '
$ws.Range("D21").Value = '11.25'
$ws.Range("E21").Value = '  +4.27%  '
$ws.Range("D22").Value = '491.30'
$ws.Range("E22").Value = '  -1.07%  '
$ws.Range("D23").Value = '0.726'
$ws.Range("E23").Value = '  -0.42%  '
$ws.Range("E24").Value = '  -0.78%  '
$ws.Range("D25").Value = '84.64'
$ws.Range("E25").Value = '  -1.18%  '
$ws.Range("D26").Value = '2.28'
$ws.Range("E26").Value = '  -2.82%  '
$ws.Range("D27").Value = '12.24'
$ws.Range("E27").Value = '  -0.95%  '
$ws.Range("D28").Value = '10.06'
$ws.Range("E28").Value = '  -1.73%  '
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("E30").Value = '  -0.71%  '
$ws.Range("D31").Value = '8.11'
$ws.Range("E31").Value = '  +1.75%  '
$ws.Range("E32").Value = '  -3.68%  '
$ws.Range("D33").Value = '31.73'
$ws.Range("E33").Value = '  -0.76%  '
$ws.Range("D34").Value = '3.900.98'
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("D35").Value = '3.698.89'
$ws.Range("E35").Value = '  +0.32%  '
$ws.Range("E36").Value = '  -0.34%  '
$ws.Range("E37").Value = '  +6.11%  '
$ws.Range("D38").Value = '5.94'
$ws.Range("E38").Value = '  +1.38%  '
$ws.Range("E39").Value = '  -0.24%  '
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("D41").Value = '3.12'
$ws.Range("E41").Value = '  +8.23%  '
$ws.Range("D42").Value = '0.326'
$ws.Range("E42").Value = '  +0.40%  '
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").Value = '427.72'
$ws.Range("E43").Value = '  -4.02%  '
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").Value = '1.99'
$ws.Range("E44").Value = '  +0.70%  '
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").Value = '48.51'
$ws.Range("E45").Value = '  -0.91%  '
$ws.Range("B46").Value = 'Cosmos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D46").Value = '8.43'
$ws.Range("E46").Value = '  -0.91%  '
$ws.Range("B47").Value = 'USDe'
$ws.Range("C47").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D47").Value = '1.00'
$ws.Range("E47").Value = '  -0.01%  '
$ws.Range("D48").Value = '142.91'
$ws.Range("E48").Value = '  +0.55%  '
$ws.Range("D49").Value = '40.24'
$ws.Range("E49").Value = '  -0.93%  '
$ws.Range("D50").Value = '1.32'
$ws.Range("E50").Value = '  +7.48%  '
$ws.Range("D51").Value = '2.798.32'
$ws.Range("E51").Value = '  -1.72%  '

# Restore default "Normal" style on column D (data rows only) so no stray number-format
# style lingers on cells (keeps styles.xml equivalent to the original, unstyled D-column
# data cells, and leaves the D1 header's bold style untouched).
$ws.Range("D2:D51").Style = "Normal"
